$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append new resale-number observation for 2024-01-06 09:28:08 (Saturday) as row 24
$row = 24

# Columns A (Date) and D (Week) must stay as plain text like the rest of the
# table ("2024-01-06", "00") instead of being auto-converted to a date serial
# / number by Excel's type inference, so force a text number format first.
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 4).NumberFormat = "@"

$ws.Cells.Item($row, 1).Value = "2024-01-06"
$ws.Cells.Item($row, 2).Value = "09:28:08"
$ws.Cells.Item($row, 3).Value = "Saturday"
$ws.Cells.Item($row, 4).Value = "00"

# Match the (unstyled) look of the other data rows now that the text is in
# place, instead of leaving the temporary "Text" number format applied.
$ws.Cells.Item($row, 1).Style = $ws.Cells.Item(2, 1).Style
$ws.Cells.Item($row, 4).Style = $ws.Cells.Item(2, 4).Style

$ws.Cells.Item($row, 5).Value = 140212
$ws.Cells.Item($row, 6).Value = 143021
$ws.Cells.Item($row, 7).Value = 171719
$ws.Cells.Item($row, 8).Value = 147162
$ws.Cells.Item($row, 9).Value = -1
$ws.Cells.Item($row, 10).Value = 117943
$ws.Cells.Item($row, 11).Value = 224398
$ws.Cells.Item($row, 12).Value = 248739
$ws.Cells.Item($row, 13).Value = 184901
$ws.Cells.Item($row, 14).Value = 110234
$ws.Cells.Item($row, 15).Value = 40539
$ws.Cells.Item($row, 16).Value = 30816
$ws.Cells.Item($row, 17).Value = 72425
$ws.Cells.Item($row, 18).Value = -1
$ws.Cells.Item($row, 19).Value = 41538
$ws.Cells.Item($row, 20).Value = -1
